# Apply "updated main GSC export data" edit to the Chart sheet.
# The date window (column A) shifts forward by one day: the oldest date
# (2025-10-25) drops off and a new date (2026-01-22) is appended; every
# row's Validation count (column C) likewise shifts up by one row, with a
# brand-new count for the newly appended date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$dates = @("2025-10-26","2025-10-27","2025-10-28","2025-10-29","2025-10-30","2025-10-31","2025-11-01","2025-11-02","2025-11-03","2025-11-04","2025-11-05","2025-11-06","2025-11-07","2025-11-08","2025-11-09","2025-11-10","2025-11-11","2025-11-12","2025-11-13","2025-11-14","2025-11-15","2025-11-16","2025-11-17","2025-11-18","2025-11-19","2025-11-20","2025-11-21","2025-11-22","2025-11-23","2025-11-24","2025-11-25","2025-11-26","2025-11-27","2025-11-28","2025-11-29","2025-11-30","2025-12-01","2025-12-02","2025-12-03","2025-12-04","2025-12-05","2025-12-06","2025-12-07","2025-12-08","2025-12-09","2025-12-10","2025-12-11","2025-12-12","2025-12-13","2025-12-14","2025-12-15","2025-12-16","2025-12-17","2025-12-18","2025-12-19","2025-12-20","2025-12-21","2025-12-22","2025-12-23","2025-12-24","2025-12-25","2025-12-26","2025-12-27","2025-12-28","2025-12-29","2025-12-30","2025-12-31","2026-01-01","2026-01-02","2026-01-03","2026-01-04","2026-01-05","2026-01-06","2026-01-07","2026-01-08","2026-01-09","2026-01-10","2026-01-11","2026-01-12","2026-01-13","2026-01-14","2026-01-15","2026-01-16","2026-01-17","2026-01-18","2026-01-19","2026-01-20","2026-01-21","2026-01-22")
$counts = @(90.0,83.0,90.0,93.0,92.0,101.0,112.0,115.0,108.0,105.0,101.0,95.0,87.0,82.0,76.0,54.0,47.0,43.0,41.0,38.0,35.0,31.0,29.0,26.0,26.0,25.0,26.0,26.0,25.0,25.0,27.0,27.0,27.0,27.0,27.0,27.0,27.0,27.0,26.0,25.0,25.0,25.0,26.0,26.0,27.0,29.0,29.0,30.0,30.0,31.0,31.0,32.0,31.0,31.0,32.0,32.0,32.0,32.0,30.0,31.0,32.0,32.0,28.0,28.0,28.0,28.0,30.0,29.0,28.0,28.0,27.0,27.0,27.0,27.0,27.0,27.0,26.0,26.0,26.0,26.0,26.0,26.0,26.0,25.0,25.0,26.0,25.0,24.0,23.0)

$firstRow = 2
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $firstRow + $i
    $cell = $ws.Cells.Item($row, 1)
    # Force text so the yyyy-MM-dd string isn't auto-converted to a date
    # serial, then drop the format change so the cell keeps its original
    # (default) style.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $cell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $counts[$i]
}
